$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits right
#    after the H1 title paragraph.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Meta description")) {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Insert a new bold paragraph ("Play Forge of Fortunes for Free -
#    Exciting Respins & Mining Theme") right before the final "Prompt for
#    DALLE" paragraph.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Forge of Fortunes for Free - Exciting Respins &amp; Mining Theme</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($xml)

# The inserted text landed in the same paragraph as "Prompt for DALLE" -
# split it into its own paragraph right before that text.
$mergedPara = $d.Paragraphs.Item($count)
$splitRange = $d.Range($mergedPara.Range.Start, $mergedPara.Range.End)
$splitRange.Find.Execute("Prompt for DALLE", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitRange.Collapse(1)
$splitRange.InsertParagraphBefore()

# ---------------------------------------------------------------------------
# 3) Replace the (still-italic) "Prompt for DALLE: ..." text with the new
#    meta-description copy.
# ---------------------------------------------------------------------------
$oldText = "Prompt for DALLE: Create an eye-catching feature image for the game Forge of Fortunes that fits its unique theme. The image should be in a cartoon style and feature a happy Maya warrior sporting glasses. The warrior should be surrounded by gold nuggets, coal, and slag to highlight the game's theme. Make sure to use bright colors and unique design elements to capture the attention of potential players. The image should be optimized for use on websites and social media platforms to promote the game in a visually appealing way."
$newText = "Experience simple gameplay with the Respins feature and win up to 2,500x in Forge of Fortunes. Play free and enjoy a gold mining theme."
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

Write-Output "done"
